$wb = $excel.ActiveWorkbook

# --- Ausgaben sheet: two new rows (19 and 20) ---
$wsAusgaben = $wb.Worksheets.Item("Ausgaben")

# Row 19
$wsAusgaben.Range("A19").Value = "Verleiher"
$wsAusgaben.Range("B19").Value = 45340
$wsAusgaben.Range("C19").Value = "Film: Planet Hora"
$wsAusgaben.Range("D19").Value = 45341
$wsAusgaben.Range("D19").NumberFormat = "m/d/yy"
$wsAusgaben.Range("E19").Value = 1000
$wsAusgaben.Range("E19").NumberFormat = """CHF""\ #,##0.00"
$wsAusgaben.Range("F19").Value = "Theater Hora"
$wsAusgaben.Range("H19").NumberFormat = "@"
$wsAusgaben.Range("I19").NumberFormat = "@"

# Row 20
$wsAusgaben.Range("A20").Value = "Eventausgaben"
$wsAusgaben.Range("B20").Value = 45340
$wsAusgaben.Range("C20").Value = "Plakate & Flyer Planet Hora"
$wsAusgaben.Range("D20").Value = 45340
$wsAusgaben.Range("D20").NumberFormat = "m/d/yy"
$wsAusgaben.Range("E20").Value = 400
$wsAusgaben.Range("E20").NumberFormat = """CHF""\ #,##0.00"
$wsAusgaben.Range("F20").Value = "L'equippe visuelle"
$wsAusgaben.Range("H20").NumberFormat = "@"
$wsAusgaben.Range("I20").NumberFormat = "@"

# Ausgaben table grows to include the new rows
[void]$wb.Worksheets.Item("Ausgaben").ListObjects.Item("Table16").Resize($wsAusgaben.Range("A1:I20"))

# Select cell E20 to match final view state
[void]$wsAusgaben.Range("E20").Select()

# --- Einnahmen sheet: one new row (3) ---
$wsEinnahmen = $wb.Worksheets.Item("Einnahmen")

$wsEinnahmen.Range("A3").Value = "Vermietung"
$wsEinnahmen.Range("E3").Value = "Stiftung Lebenshilfe"
$wsEinnahmen.Range("F3").Value = "Reinach AG"
$wsEinnahmen.Range("B3").Value = "Beitrag Stiftung Lebenshilfe"
$wsEinnahmen.Range("C3").Value = 45340
$wsEinnahmen.Range("D3").Value = 1000

# Select cell B3 and mark Einnahmen as the active/tab-selected sheet
[void]$wsEinnahmen.Activate()
[void]$wsEinnahmen.Range("B3").Select()
